# Applies the diff: reorders/expands the "27.03.2020" journal entry by
# adding two new bullet points under "Aufgetretene Probleme / Misserfolge"
# and one new paragraph under "Reflexion", moving the _GoBack bookmark
# along the way.
$d = $word.ActiveDocument

# --- Helper: find paragraph index (1-based) whose text contains a marker ---
function Find-ParaIndex($doc, [string]$marker) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like "*$marker*") {
            return $idx
        }
    }
    return -1
}

# 1. Remove the original _GoBack bookmark; it gets re-created further up
#    the document once the new bullet paragraph exists.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the "Unit Tests aufsetzen..." bullet paragraph (last bullet of
#    the "Aufgetretene Probleme / Misserfolge" section, 27.03.2020 entry).
$bulletIdx = Find-ParaIndex $d "Unit Tests aufsetzen funktionierte nicht so wie es sollte."

# 3. Insert the first new bullet paragraph's first sentence right after it.
$d.Paragraphs.Item($bulletIdx).Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($bulletIdx + 1)
$p1.Range.Text = "Führe ich die Unit Tests einzeln aus funktionieren alle, führe ich jedoch alle nacheinander aus so wird bei welchen die Datenbank nicht verbunden."

# 3b. Add the second sentence as its own paragraph, then merge the paragraph
#     break away so the two sentences end up as two separate runs inside a
#     single paragraph (mirrors how the source document is structured).
$d.Paragraphs.Item($bulletIdx + 1).Range.InsertParagraphAfter()
$p1b = $d.Paragraphs.Item($bulletIdx + 2)
$p1b.Range.Text = " Dies ist nur bei Tests mit Starlette so."
$p1again = $d.Paragraphs.Item($bulletIdx + 1)
$mergeRange = $d.Range($p1again.Range.End - 1, $p1again.Range.End)
$mergeRange.Delete()

# 4. Insert the second new bullet paragraph after the first one.
$d.Paragraphs.Item($bulletIdx + 1).Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($bulletIdx + 2)
$p2.Range.Text = "Es fehlte eine Spalte für die Stockwerke Tabelle"

# 5. Re-create the _GoBack bookmark at the very end of this paragraph's text
#    (immediately before its paragraph mark), matching the source layout.
#    A zero-length bookmark at that exact boundary can't be added directly,
#    so a temporary placeholder character is inserted after the text, the
#    bookmark is added just before it, and the placeholder is removed again.
$p2 = $d.Paragraphs.Item($bulletIdx + 2)
$endPos = $p2.Range.End - 1
$placeholderRange = $d.Range($endPos, $endPos)
$placeholderRange.InsertAfter("X")
$p2 = $d.Paragraphs.Item($bulletIdx + 2)
$bmPos = $d.Range($p2.Range.End - 2, $p2.Range.End - 2)
$d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null
$p2 = $d.Paragraphs.Item($bulletIdx + 2)
$placeholderDelRange = $d.Range($p2.Range.End - 2, $p2.Range.End - 1)
$placeholderDelRange.Delete()

# 6. Locate the "Von allen Sachen..." reflection paragraph and add the new
#    paragraph about Flüchtigkeitsfehler right after it.
$reflectionIdx = Find-ParaIndex $d "Von allen Sachen welche Probleme erstellen könnten war das Aufsetzen von Unit Tests"
$d.Paragraphs.Item($reflectionIdx).Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($reflectionIdx + 1)
$p3.Range.Text = "Ich hatte heute sehr viele Flüchtigkeitsfehler welche einfach verhindert werden hätten könnte. Ich hoffe diese kann ich am Montag mit neuer Kraft vermeiden. "
